$wb = $excel.ActiveWorkbook

# 1) Status text change: "Ready for handoff" -> "In Translation" (all sheets, all occurrences)
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# 2) Column width shrink that accompanies the shorter status text:
#    Overview sheet columns E (zh-cn) and F (de-de), and the per-locale
#    sheets' "Status" column (C) all go from ~17.22 chars to ~13.41 chars.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.43
$wsOverview.Columns.Item(6).ColumnWidth = 12.43

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.43

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.43
